$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.955.25"
$ws.Range("E2").Value = "  -4.22%  "

$ws.Range("D3").Value = "2.327.87"
$ws.Range("E3").Value = "  -5.86%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("E5").Value = "  -4.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "84.61"
$ws.Range("E6").Value = "  -8.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.530"
$ws.Range("E7").Value = "  -3.66%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.483"
$ws.Range("E9").Value = "  -5.08%  "

$ws.Range("E10").Value = "  -4.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.01"
$ws.Range("E11").Value = "  -8.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.110"
$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").Value = "2.686.64"
$ws.Range("E13").Value = "  -5.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.41"
$ws.Range("E14").Value = "  -6.75%  "

$ws.Range("E15").Value = "  -5.48%  "

$ws.Range("D16").Value = "2.328.93"
$ws.Range("E16").Value = "  -5.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.750"
$ws.Range("E17").Value = "  -4.96%  "

$ws.Range("D18").Value = "39.927.81"
$ws.Range("E18").Value = "  -4.08%  "

$ws.Range("E19").Value = "  -4.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.07"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.58"
$ws.Range("E21").Value = "  -5.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.64"
$ws.Range("E22").Value = "  -5.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.49"
$ws.Range("E23").Value = "  -1.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.55"
$ws.Range("E24").Value = "  -7.27%  "

$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.81"
$ws.Range("E26").Value = "  -6.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.37"
$ws.Range("E27").Value = "  -5.92%  "

$ws.Range("E28").Value = "  -1.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.25"
$ws.Range("E29").Value = "  -4.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.43"
$ws.Range("E30").Value = "  -2.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.00"
$ws.Range("E31").Value = "  -1.95%  "

$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.10"
$ws.Range("E33").Value = "  -5.99%  "

$ws.Range("E34").Value = "  -4.31%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0722"
$ws.Range("E35").Value = "  -5.33%  "

$ws.Range("E36").Value = "  -2.13%  "

$ws.Range("E37").Value = "  -3.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.75"
$ws.Range("E38").Value = "  -5.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.75"
$ws.Range("E39").Value = "  -7.87%  "

$ws.Range("E40").Value = "  -7.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.81"
$ws.Range("E41").Value = "  -4.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.27"
$ws.Range("E42").Value = "  -6.13%  "

$ws.Range("D43").Value = "1.942.20"
$ws.Range("E43").Value = "  -2.98%  "

$ws.Range("E44").Value = "  -5.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.53"
$ws.Range("E45").Value = "  -5.76%  "

$ws.Range("E46").Value = "  -1.77%  "

$ws.Range("E47").Value = "  -9.33%  "

$ws.Range("D48").Value = "2.556.79"
$ws.Range("E48").Value = "  -6.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.69"
$ws.Range("E49").Value = "  -4.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.60"
$ws.Range("E50").Value = "  -6.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.34"
$ws.Range("E51").Value = "  -3.29%  "
